# OpenTBS 1.9.1-beta-2014-07-22 : credit features
#
# The canonical diff for this commit only touches the 3-D bar chart on
# slide 3 (the "Graphique" shape): the two internal axis identifiers
# (<c:axId>/<c:crossAx> in ppt/charts/chart1.xml) are renumbered from
# 95843456 / 95844992 to 61990016 / 61991552, consistently, everywhere
# they are referenced (bar3DChart axId list, catAx/valAx axId and the
# matching crossAx backlinks). No visible chart content, formatting,
# series data, or axis scaling changes.
#
# Axis ids are an internal bookkeeping value assigned by PowerPoint
# when a chart is (re)built; there is no PowerPoint object-model
# property that exposes them for editing (Axis has no settable
# AxisID/Id property) - this host environment explicitly rejects a
# write to it ("... is not supported for chart parts in this
# environment"). We still try it defensively (future-proof / no-op if
# unsupported) instead of touching anything else, since there is no
# other observable-content change requested by this commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$chartShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasChart) {
        $chartShape = $candidate
    }
}

if ($chartShape -ne $null) {
    $chart = $chartShape.Chart

    $catAxis = $chart.Axes(1)  # xlCategory -> <c:catAx>
    $valAxis = $chart.Axes(2)  # xlValue    -> <c:valAx>

    try { $catAxis.AxisID = 61990016 } catch { }
    try { $valAxis.AxisID = 61991552 } catch { }
}
